$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all changed Price (D) cells to Text format first to prevent Excel
# from auto-converting numeric-looking strings into floating point numbers
# (which would introduce precision artifacts / scientific notation).
$textCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D12","D14","D15","D16","D18","D19","D20","D21","D22","D23","D24","D25","D26","D29","D30","D31","D32","D34","D35","D37","D38","D40","D41","D42","D46","D47","D48","D49","D50","D51","D27","D28")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "54.769.40"
$ws.Range("E2").Value = "  +1.22%  "

$ws.Range("D3").Value = "2.302.48"
$ws.Range("E3").Value = "  +0.74%  "

$ws.Range("D4").Value = "0.995"
$ws.Range("E4").Value = "  -0.69%  "

$ws.Range("D5").Value = "498.66"
$ws.Range("E5").Value = "  +1.23%  "

$ws.Range("D6").Value = "128.87"
$ws.Range("E6").Value = "  +0.82%  "

$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  -0.48%  "

$ws.Range("D8").Value = "0.531"
$ws.Range("E8").Value = "  +0.79%  "

$ws.Range("D9").Value = "2.299.90"
$ws.Range("E9").Value = "  +0.33%  "

$ws.Range("E10").Value = "  +2.00%  "

$ws.Range("E11").Value = "  +2.39%  "

$ws.Range("D12").Value = "0.324"
$ws.Range("E12").Value = "  +2.24%  "

$ws.Range("E13").Value = "  -1.75%  "

$ws.Range("D14").Value = "2.689.79"
$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("D15").Value = "21.81"
$ws.Range("E15").Value = "  +2.37%  "

$ws.Range("D16").Value = "54.371.84"
$ws.Range("E16").Value = "  +0.49%  "

$ws.Range("E17").Value = "  +0.45%  "

$ws.Range("D18").Value = "2.309.51"
$ws.Range("E18").Value = "  +0.48%  "

$ws.Range("D19").Value = "10.11"
$ws.Range("E19").Value = "  +4.67%  "

$ws.Range("D20").Value = "4.12"
$ws.Range("E20").Value = "  +2.83%  "

$ws.Range("D21").Value = "307.96"
$ws.Range("E21").Value = "  +1.58%  "

$ws.Range("D22").Value = "6.46"
$ws.Range("E22").Value = "  +4.83%  "

$ws.Range("D23").Value = "0.995"
$ws.Range("E23").Value = "  -0.36%  "

$ws.Range("D24").Value = "5.36"
$ws.Range("E24").Value = "  -1.69%  "

$ws.Range("D25").Value = "62.80"
$ws.Range("E25").Value = "  -1.84%  "

$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.24%  "

$ws.Range("D29").Value = "2.407.86"
$ws.Range("E29").Value = "  +0.29%  "

$ws.Range("D30").Value = "7.16"
$ws.Range("E30").Value = "  +1.02%  "

$ws.Range("D31").Value = "170.04"

$ws.Range("D32").Value = "0.0₃0698"
$ws.Range("E32").Value = "  -0.06%  "

$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("D34").Value = "5.94"
$ws.Range("E34").Value = "  +2.67%  "

$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  -0.16%  "

$ws.Range("E36").Value = "  +1.20%  "

$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").Value = "17.69"
$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("E39").Value = "  +2.98%  "

$ws.Range("D40").Value = "0.872"
$ws.Range("E40").Value = "  +2.94%  "

$ws.Range("D41").Value = "3.69"
$ws.Range("E41").Value = "  +1.66%  "

$ws.Range("D42").Value = "35.49"
$ws.Range("E42").Value = "  -0.83%  "

$ws.Range("E43").Value = "  +2.75%  "

$ws.Range("E44").Value = "  +2.17%  "

$ws.Range("E45").Value = "  +0.92%  "

$ws.Range("D46").Value = "128.45"
$ws.Range("E46").Value = "  +3.94%  "

$ws.Range("D47").Value = "4.93"
$ws.Range("E47").Value = "  +3.59%  "

$ws.Range("D48").Value = "0.0894"
$ws.Range("E48").Value = "  +1.34%  "

$ws.Range("D49").Value = "0.551"
$ws.Range("E49").Value = "  +0.92%  "

$ws.Range("D50").Value = "244.90"
$ws.Range("E50").Value = "  +2.98%  "

$ws.Range("D51").Value = "0.0487"
$ws.Range("E51").Value = "  +2.70%  "

# Row 27 and 28: Kaspa/Polygon swap
$ws.Range("B27").Value = "Polygon"
$ws.Range("C27").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D27").Value = "0.375"
$ws.Range("E27").Value = "  +2.39%  "

$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.152"
$ws.Range("E28").Value = "  +6.09%  "